$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common (constant across all rows in this dataset)
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$rubroId = 100103
$rubro = "Frutos de hueso (carozo)"
$especieId = 100103001
$especie = "Cereza"
$unidad = "`$/bandeja 10 kilos"
$codigo = 10

function Set-Row($r, $fecha, $variedad, $calidad, $calibre, $precioMin, $precioProm, $precioMax, $regionOrigen, $precioKilo) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $rubroId
    $ws.Cells.Item($r, 8).Value = $rubro
    $ws.Cells.Item($r, 9).Value = $especieId
    $ws.Cells.Item($r, 10).Value = $especie
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $calibre
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioProm
    $ws.Cells.Item($r, 16).Value = $precioMax
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $regionOrigen
    $ws.Cells.Item($r, 19).Value = $precioKilo
    $ws.Cells.Item($r, 20).Value = $codigo
}

# Rows 530-532: updated with new "Brooks" price records (dated 2022-11-25 / serial 44890)
Set-Row 530 44890 "Brooks" "Especial" 67 15000 15000 15000 "Región de O'Higgins" 1500
Set-Row 531 44890 "Brooks" "Primera" 60 12000 12000 12000 "Región de O'Higgins" 1200
Set-Row 532 44890 "Brooks" "Segunda" 50 10000 10000 10000 "Región de O'Higgins" 1000

# Rows 533-535 (new): "Sweet Heart" price records (dated 2022-11-25 / serial 44890)
Set-Row 533 44890 "Sweet Heart" "Especial" 60 15000 15000 15000 "Región de O'Higgins" 1500
Set-Row 534 44890 "Sweet Heart" "Primera" 60 12000 12000 12000 "Región de O'Higgins" 1200
Set-Row 535 44890 "Sweet Heart" "Segunda" 50 10000 10000 10000 "Región de O'Higgins" 1000

# Rows 536-538 (new): original rows 530-532 data, shifted down unchanged
Set-Row 536 44225 "Santina" "Especial" 45 13000 13000 13000 "Región del Maule" 1300
Set-Row 537 44225 "Santina" "Primera" 56 12000 12000 12000 "Región del Maule" 1200
Set-Row 538 44236 "Rainier" "Primera" 150 10000 10000 10000 "Provincia de Curicó" 1000
